$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.502.29"
$ws.Range("E2").Value = "'  -0.79%  "
$ws.Range("D3").Value = "'1.831.75"
$ws.Range("E3").Value = "'  -0.87%  "
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("D5").Value = "'312.58"
$ws.Range("E5").Value = "'  -0.50%  "
$ws.Range("E6").Value = "'  +0.02%  "
$ws.Range("D7").Value = "'0.4291"
$ws.Range("E7").Value = "'  -0.31%  "
$ws.Range("D8").Value = "'0.3662"
$ws.Range("D9").Value = "'0.07291"
$ws.Range("E9").Value = "'  -0.63%  "
$ws.Range("D10").Value = "'0.8665"
$ws.Range("E10").Value = "'  -1.14%  "
$ws.Range("D11").Value = "'20.69"
$ws.Range("E11").Value = "'  -0.34%  "
$ws.Range("D12").Value = "'1.899.65"
$ws.Range("E12").Value = "'  +4.37%  "
$ws.Range("D13").Value = "'5.411"
$ws.Range("E13").Value = "'  +1.10%  "
$ws.Range("E14").Value = "'  -0.02%  "
$ws.Range("D15").Value = "'0.06930"
$ws.Range("E15").Value = "'  -0.22%  "
$ws.Range("E16").Value = "'  +0.25%  "
$ws.Range("D17").Value = "'80.49"
$ws.Range("E17").Value = "'  +0.77%  "
$ws.Range("D18").Value = "'0.000008924"
$ws.Range("E18").Value = "'  -0.68%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "'  +0.07%  "
$ws.Range("E20").Value = "'  +0.30%  "
$ws.Range("D21").Value = "'27.364.05"
$ws.Range("E21").Value = "'  -0.95%  "
$ws.Range("D22").Value = "'5.140"
$ws.Range("E22").Value = "'  +3.20%  "
$ws.Range("D23").Value = "'10.87"
$ws.Range("E23").Value = "'  +5.00%  "
$ws.Range("D24").Value = "'2.043.43"
$ws.Range("E24").Value = "'  -1.02%  "
$ws.Range("D25").Value = "'1.980"
$ws.Range("D26").Value = "'154.56"
$ws.Range("E26").Value = "'  -1.11%  "
$ws.Range("D27").Value = "'18.95"
$ws.Range("E27").Value = "'  +1.55%  "
$ws.Range("D28").Value = "'5.141"
$ws.Range("E28").Value = "'  -2.35%  "
$ws.Range("D29").Value = "'114.43"
$ws.Range("E29").Value = "'  -4.55%  "
$ws.Range("D30").Value = "'1.839"
$ws.Range("E30").Value = "'  -2.04%  "
$ws.Range("D31").Value = "'0.08863"
$ws.Range("E31").Value = "'  -0.42%  "
$ws.Range("D32").Value = "'0.7551"
$ws.Range("E32").Value = "'  -0.08%  "
$ws.Range("E33").Value = "'  +0.66%  "
$ws.Range("D34").Value = "'4.547"
$ws.Range("E34").Value = "'  -0.17%  "
$ws.Range("D35").Value = "'1.138"
$ws.Range("E35").Value = "'  +1.00%  "
$ws.Range("E36").Value = "'  +0.08%  "
$ws.Range("D37").Value = "'1.088"
$ws.Range("E37").Value = "'  -1.72%  "
$ws.Range("D38").Value = "'0.05335"
$ws.Range("E38").Value = "'  -1.78%  "
$ws.Range("E39").Value = "'  +0.28%  "
$ws.Range("D40").Value = "'2.799"
$ws.Range("D41").Value = "'0.5098"
$ws.Range("E41").Value = "'  +0.09%  "
$ws.Range("D42").Value = "'0.1668"
$ws.Range("E42").Value = "'  +0.26%  "
$ws.Range("D43").Value = "'6.576"
$ws.Range("E43").Value = "'  +0.01%  "
$ws.Range("D44").Value = "'8.384"
$ws.Range("E44").Value = "'  -0.04%  "
$ws.Range("D45").Value = "'10.52"
$ws.Range("E45").Value = "'  +0.95%  "
$ws.Range("D46").Value = "'106.24"
$ws.Range("E46").Value = "'  +0.76%  "
$ws.Range("D47").Value = "'0.06503"
$ws.Range("E47").Value = "'  -0.60%  "
$ws.Range("E48").Value = "'  +0.49%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "'  +0.03%  "
$ws.Range("D50").Value = "'1.620"
$ws.Range("E50").Value = "'  -0.73%  "
$ws.Range("D51").Value = "'64.08"
$ws.Range("E51").Value = "'  -0.65%  "
